# Update the build-version timestamp that appears throughout the workbook.
$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$wb = $excel.ActiveWorkbook

# --- Sheet "About" ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: mines - January 30 (built on " + $newStamp + ")"

$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Chaili Coal Mine, China, M1829, version 'mines - January 30 (built on " + $newStamp + ")'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- Sheet "Boundaries and methane sources" ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$newBuildVersion = "mines - January 30 (built on " + $newStamp + ")"

$row = 2
while ($row -le 7) {
    $wsData.Range("S" + $row).Value = $newBuildVersion
    $row = $row + 1
}
